# Update the "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1130
    3  = 841
    4  = 277
    8  = 2064
    9  = 7700
    10 = 917
    11 = 436
    12 = 373
    13 = 151
    14 = 415
    15 = 160
    16 = 7878
    17 = 318
    18 = 1371
    22 = 165
    23 = 318
    24 = 154
    25 = 165
    28 = 25
    30 = 1145
    31 = 57
    32 = 94
    33 = 66
    34 = 82
    36 = 78
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
